$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 286; existing rows 286:316 shift down to 287:317.
$ws.Rows("286:286").Insert()

# Populate the newly inserted row 286 with its data.
$ws.Range("A286").Value2 = 3
$ws.Range("B286").Value = "Femacal de La Calera"
$ws.Range("C286").Value = "Coquimbo"
$ws.Range("D286").Value2 = 45212
$ws.Range("D286").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E286").Value2 = 5
$ws.Range("F286").Value2 = 100112026
$ws.Range("G286").Value = "Haba"
$ws.Range("H286").Value = "Sin especificar"
$ws.Range("I286").Value = "Primera"
$ws.Range("J286").Value2 = 40
$ws.Range("K286").Value2 = 13000
$ws.Range("L286").Value2 = 13000
$ws.Range("M286").Value2 = 13000
$ws.Range("N286").Value = "$/saco 25 kilos"
$ws.Range("O286").Value = "Provincia de Quillota"
$ws.Range("P286").Value2 = 520
$ws.Range("Q286").Value2 = 25
$ws.Range("R286").Value = "Hortaliza"
